# Refresh market-price-derived columns (currentAveragePrice / NQ / HQ,
# LevePriceNQ/HQ, LeveProfitNQ/HQ) for leves whose Universalis price data
# changed since the last scheduled pull. Values only; no formulas/styles.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 41
$ws.Range("H41").Value = 704.3889
$ws.Range("J41").Value = 670.9091
$ws.Range("L41").Value = 670.9091
$ws.Range("N41").Value = -1550.9091
# Row 74
$ws.Range("H74").Value = 6254925
$ws.Range("I74").Value = 33336000
$ws.Range("J74").Value = 5446.154
$ws.Range("K74").Value = 33336000
$ws.Range("L74").Value = 5446.154
$ws.Range("M74").Value = -33335064
$ws.Range("N74").Value = -7318.154
# Row 77
$ws.Range("H77").Value = 6254925
$ws.Range("I77").Value = 33336000
$ws.Range("J77").Value = 5446.154
$ws.Range("K77").Value = 166680000
$ws.Range("L77").Value = 27230.77
$ws.Range("M77").Value = -166675320
$ws.Range("N77").Value = -36590.77
# Row 132
$ws.Range("H132").Value = 12471882
$ws.Range("I132").Value = 15875772
$ws.Range("J132").Value = 558268.5
$ws.Range("K132").Value = 47627316
$ws.Range("L132").Value = 1674805.5
$ws.Range("M132").Value = -47624786
$ws.Range("N132").Value = -1679865.5
# Row 137
$ws.Range("H137").Value = 1965.6227
$ws.Range("I137").Value = 1018.55
$ws.Range("J137").Value = 4879.6924
$ws.Range("K137").Value = 3055.65
$ws.Range("L137").Value = 14639.0772
$ws.Range("M137").Value = -505.6499999999996
$ws.Range("N137").Value = -19739.0772
# Row 138
$ws.Range("H138").Value = 3241.1
$ws.Range("I138").Value = 843.67645
$ws.Range("J138").Value = 4476.136
$ws.Range("K138").Value = 2531.02935
$ws.Range("L138").Value = 13428.408
$ws.Range("M138").Value = 2608.97065
$ws.Range("N138").Value = -23708.408
# Row 141
$ws.Range("H141").Value = 5583.1914
$ws.Range("I141").Value = 6195.2563
$ws.Range("K141").Value = 18585.7689
$ws.Range("M141").Value = -13405.7689

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5689.2856
$ws.Range("I32").Value = 4174.9585
$ws.Range("J32").Value = 10535.134
$ws.Range("K32").Value = 4174.9585
$ws.Range("L32").Value = 10535.134
$ws.Range("M32").Value = -3887.9585
$ws.Range("N32").Value = -11109.134
# Row 38
$ws.Range("H38").Value = 9729.75
$ws.Range("I38").Value = 4639.6665
$ws.Range("K38").Value = 4639.6665
$ws.Range("M38").Value = -4172.6665
# Row 54
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
# Row 61
$ws.Range("H61").Value = 1088.1034
$ws.Range("I61").Value = 951.5625
$ws.Range("J61").Value = 1256.1538
$ws.Range("K61").Value = 951.5625
$ws.Range("L61").Value = 1256.1538
$ws.Range("M61").Value = -739.5625
$ws.Range("N61").Value = -1680.1538
# Row 122
$ws.Range("H122").Value = 2416.1738
$ws.Range("I122").Value = 1298.7142
$ws.Range("J122").Value = 4154.4443
$ws.Range("K122").Value = 3896.1426
$ws.Range("L122").Value = 12463.3329
$ws.Range("M122").Value = -1446.1426
$ws.Range("N122").Value = -17363.3329
# Row 132
$ws.Range("H132").Value = 1709.6038
$ws.Range("I132").Value = 1145.8478
$ws.Range("J132").Value = 5414.2856
$ws.Range("K132").Value = 3437.5434
$ws.Range("L132").Value = 16242.8568
$ws.Range("M132").Value = -907.5434
$ws.Range("N132").Value = -21302.8568
# Row 136
$ws.Range("H136").Value = 1088.1034
$ws.Range("I136").Value = 951.5625
$ws.Range("J136").Value = 1256.1538
$ws.Range("K136").Value = 2854.6875
$ws.Range("L136").Value = 3768.4614
$ws.Range("M136").Value = -304.6875
$ws.Range("N136").Value = -8868.4614

$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 625.25
$ws.Range("I22").Value = 625.25
$ws.Range("K22").Value = 625.25
$ws.Range("M22").Value = -452.25
# Row 38
$ws.Range("H38").Value = 18017.5
$ws.Range("J38").Value = 18017.5
$ws.Range("L38").Value = 18017.5
$ws.Range("N38").Value = -18849.5
# Row 44
$ws.Range("H44").Value = 24021.8
$ws.Range("J44").Value = 24021.8
$ws.Range("L44").Value = 24021.8
$ws.Range("N44").Value = -25015.8
# Row 49
$ws.Range("H49").Value = 11999
$ws.Range("J49").Value = 11999
$ws.Range("L49").Value = 11999
$ws.Range("N49").Value = -12477
# Row 75
$ws.Range("H75").Value = 9675.5
$ws.Range("I75").Value = 4951
$ws.Range("J75").Value = 14400
$ws.Range("K75").Value = 4951
$ws.Range("L75").Value = 14400
$ws.Range("M75").Value = -4015
$ws.Range("N75").Value = -16272
# Row 76
$ws.Range("H76").Value = 20000
$ws.Range("J76").Value = 20000
$ws.Range("L76").Value = 20000
$ws.Range("N76").Value = -20630
# Row 78
$ws.Range("H78").Value = 9675.5
$ws.Range("I78").Value = 4951
$ws.Range("J78").Value = 14400
$ws.Range("K78").Value = 14853
$ws.Range("L78").Value = 43200
$ws.Range("M78").Value = -10173
$ws.Range("N78").Value = -52560
# Row 79
$ws.Range("H79").Value = 20000
$ws.Range("J79").Value = 20000
$ws.Range("L79").Value = 20000
$ws.Range("N79").Value = -22184

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 10640707
$ws.Range("I31").Value = 1374.8966
$ws.Range("J31").Value = 27781854
$ws.Range("K31").Value = 1374.8966
$ws.Range("L31").Value = 27781854
$ws.Range("M31").Value = -1079.8966
$ws.Range("N31").Value = -27782444
# Row 34
$ws.Range("H34").Value = 10640707
$ws.Range("I34").Value = 1374.8966
$ws.Range("J34").Value = 27781854
$ws.Range("K34").Value = 1374.8966
$ws.Range("L34").Value = 27781854
$ws.Range("M34").Value = -1172.8966
$ws.Range("N34").Value = -27782258
# Row 132
$ws.Range("H132").Value = 1998.6552
$ws.Range("I132").Value = 1106.6666
$ws.Range("K132").Value = 3319.9998
$ws.Range("M132").Value = -789.9998000000001
# Row 134
$ws.Range("H134").Value = 3495.86
$ws.Range("I134").Value = 3898.3125
$ws.Range("K134").Value = 11694.9375
$ws.Range("M134").Value = -9159.9375

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1229.75
$ws.Range("I5").Value = 434.6087
$ws.Range("K5").Value = 1303.8261
$ws.Range("M5").Value = -1191.8261
# Row 41
$ws.Range("H41").Value = 433.33334
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
# Row 135
$ws.Range("H135").Value = 1229.75
$ws.Range("I135").Value = 434.6087
$ws.Range("K135").Value = 3911.4783
$ws.Range("M135").Value = -1376.4783

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 2371.139
$ws.Range("I132").Value = 1369.6666
$ws.Range("J132").Value = 5375.5557
$ws.Range("K132").Value = 4108.9998
$ws.Range("L132").Value = 16126.6671
$ws.Range("M132").Value = -1578.9998
$ws.Range("N132").Value = -21186.6671

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 2366.6667
$ws.Range("I100").Value = 2050
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2050
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -1509
$ws.Range("N100").Value = -4082
# Row 141
$ws.Range("H141").Value = 41177.69
$ws.Range("J141").Value = 41177.69
$ws.Range("L141").Value = 41177.69
$ws.Range("N141").Value = -51537.69

$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 5954036.5
$ws.Range("I132").Value = 1215.6428
$ws.Range("K132").Value = 3646.9284
$ws.Range("M132").Value = -1116.9284
